# End of day commit: advance the date/timestamp by one day for all rate
# rows and update the ANG rate value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date (column D), time_last_updated (column E), for rows 2-6
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 4).Value = 44534
    $ws.Cells.Item($row, 5).Value = 1638576002
}

# Update the ANG rate (row 5, column F)
$ws.Cells.Item(5, 6).Value = 489.42
